$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("2026-01-24", 9, 0, "24.01.20269"),
  @("2026-01-24", 10, 0, "24.01.202610"),
  @("2026-01-24", 11, 0.098, "24.01.202611"),
  @("2026-01-24", 12, 0.372, "24.01.202612"),
  @("2026-01-24", 13, 0.513, "24.01.202613"),
  @("2026-01-24", 14, 0.53, "24.01.202614"),
  @("2026-01-24", 15, 0.54, "24.01.202615"),
  @("2026-01-24", 16, 0.43, "24.01.202616"),
  @("2026-01-24", 17, 0.191, "24.01.202617"),
  @("2026-01-24", 18, 0.057, "24.01.202618"),
  @("2026-01-24", 19, 0, "24.01.202619"),
  @("2026-01-24", 20, 0, "24.01.202620"),
  @("2026-01-24", 21, 0, "24.01.202621"),
  @("2026-01-24", 22, 0, "24.01.202622"),
  @("2026-01-24", 23, 0, "24.01.202623"),
  @("2026-01-24", 24, 0, "24.01.202624"),
  @("2026-01-25", 1, 0, "25.01.20261"),
  @("2026-01-25", 2, 0, "25.01.20262"),
  @("2026-01-25", 3, 0, "25.01.20263"),
  @("2026-01-25", 4, 0, "25.01.20264"),
  @("2026-01-25", 5, 0, "25.01.20265"),
  @("2026-01-25", 6, 0, "25.01.20266"),
  @("2026-01-25", 7, 0, "25.01.20267"),
  @("2026-01-25", 8, 0, "25.01.20268"),
  @("2026-01-25", 9, 0, "25.01.20269"),
  @("2026-01-25", 10, 0.05, "25.01.202610"),
  @("2026-01-25", 11, 0.28, "25.01.202611"),
  @("2026-01-25", 12, 0.405, "25.01.202612"),
  @("2026-01-25", 13, 0.874, "25.01.202613"),
  @("2026-01-25", 14, 1.047, "25.01.202614"),
  @("2026-01-25", 15, 1.006, "25.01.202615"),
  @("2026-01-25", 16, 0.473, "25.01.202616"),
  @("2026-01-25", 17, 0.205, "25.01.202617"),
  @("2026-01-25", 18, 0.052, "25.01.202618"),
  @("2026-01-25", 19, 0, "25.01.202619"),
  @("2026-01-25", 20, 0, "25.01.202620"),
  @("2026-01-25", 21, 0, "25.01.202621"),
  @("2026-01-25", 22, 0, "25.01.202622"),
  @("2026-01-25", 23, 0, "25.01.202623"),
  @("2026-01-25", 24, 0, "25.01.202624"),
  @("2026-01-26", 1, 0, "26.01.20261"),
  @("2026-01-26", 2, 0, "26.01.20262"),
  @("2026-01-26", 3, 0, "26.01.20263"),
  @("2026-01-26", 4, 0, "26.01.20264"),
  @("2026-01-26", 5, 0, "26.01.20265"),
  @("2026-01-26", 6, 0, "26.01.20266"),
  @("2026-01-26", 7, 0, "26.01.20267"),
  @("2026-01-26", 8, 0, "26.01.20268"),
  @("2026-01-26", 9, 0, "26.01.20269"),
  @("2026-01-26", 10, 0.017, "26.01.202610"),
  @("2026-01-26", 11, 0.109, "26.01.202611"),
  @("2026-01-26", 12, 0.269, "26.01.202612"),
  @("2026-01-26", 13, 0.405, "26.01.202613"),
  @("2026-01-26", 14, 0.411, "26.01.202614"),
  @("2026-01-26", 15, 0.413, "26.01.202615"),
  @("2026-01-26", 16, 0.22, "26.01.202616"),
  @("2026-01-26", 17, 0.111, "26.01.202617"),
  @("2026-01-26", 18, 0.021, "26.01.202618"),
  @("2026-01-26", 19, 0, "26.01.202619"),
  @("2026-01-26", 20, 0, "26.01.202620"),
  @("2026-01-26", 21, 0, "26.01.202621"),
  @("2026-01-26", 22, 0, "26.01.202622"),
  @("2026-01-26", 23, 0, "26.01.202623"),
  @("2026-01-26", 24, 0, "26.01.202624"),
  @("2026-01-27", 1, 0, "27.01.20261"),
  @("2026-01-27", 2, 0, "27.01.20262"),
  @("2026-01-27", 3, 0, "27.01.20263"),
  @("2026-01-27", 4, 0, "27.01.20264"),
  @("2026-01-27", 5, 0, "27.01.20265"),
  @("2026-01-27", 6, 0, "27.01.20266"),
  @("2026-01-27", 7, 0, "27.01.20267"),
  @("2026-01-27", 8, 0, "27.01.20268"),
  @("2026-01-27", 9, 0, "27.01.20269"),
  @("2026-01-27", 10, 0.02, "27.01.202610"),
  @("2026-01-27", 11, 0.107, "27.01.202611"),
  @("2026-01-27", 12, 0.181, "27.01.202612"),
  @("2026-01-27", 13, 0.312, "27.01.202613"),
  @("2026-01-27", 14, 0.394, "27.01.202614"),
  @("2026-01-27", 15, 0.366, "27.01.202615"),
  @("2026-01-27", 16, 0.231, "27.01.202616"),
  @("2026-01-27", 17, 0.136, "27.01.202617"),
  @("2026-01-27", 18, 0.024, "27.01.202618"),
  @("2026-01-27", 19, 0, "27.01.202619"),
  @("2026-01-27", 20, 0, "27.01.202620"),
  @("2026-01-27", 21, 0, "27.01.202621"),
  @("2026-01-27", 22, 0, "27.01.202622"),
  @("2026-01-27", 23, 0, "27.01.202623"),
  @("2026-01-27", 24, 0, "27.01.202624"),
  @("2026-01-28", 1, 0, "28.01.20261"),
  @("2026-01-28", 2, 0, "28.01.20262"),
  @("2026-01-28", 3, 0, "28.01.20263"),
  @("2026-01-28", 4, 0, "28.01.20264"),
  @("2026-01-28", 5, 0, "28.01.20265"),
  @("2026-01-28", 6, 0, "28.01.20266"),
  @("2026-01-28", 7, 0, "28.01.20267"),
  @("2026-01-28", 8, 0, "28.01.20268"),
  @("2026-01-28", 9, 0, "28.01.20269"),
  @("2026-01-28", 10, 0.028, "28.01.202610"),
  @("2026-01-28", 11, 0.132, "28.01.202611"),
  @("2026-01-28", 12, 0.424, "28.01.202612"),
  @("2026-01-28", 13, 0.6, "28.01.202613"),
  @("2026-01-28", 14, 0.84, "28.01.202614"),
  @("2026-01-28", 15, 0.769, "28.01.202615"),
  @("2026-01-28", 16, 0.58, "28.01.202616"),
  @("2026-01-28", 17, 0.286, "28.01.202617"),
  @("2026-01-28", 18, 0.059, "28.01.202618"),
  @("2026-01-28", 19, 0, "28.01.202619"),
  @("2026-01-28", 20, 0, "28.01.202620"),
  @("2026-01-28", 21, 0, "28.01.202621"),
  @("2026-01-28", 22, 0, "28.01.202622"),
  @("2026-01-28", 23, 0, "28.01.202623"),
  @("2026-01-28", 24, 0, "28.01.202624"),
  @("2026-01-29", 1, 0, "29.01.20261"),
  @("2026-01-29", 2, 0, "29.01.20262"),
  @("2026-01-29", 3, 0, "29.01.20263"),
  @("2026-01-29", 4, 0, "29.01.20264"),
  @("2026-01-29", 5, 0, "29.01.20265"),
  @("2026-01-29", 6, 0, "29.01.20266"),
  @("2026-01-29", 7, 0, "29.01.20267"),
  @("2026-01-29", 8, 0, "29.01.20268"),
  @("2026-01-29", 9, 0, "29.01.20269"),
  @("2026-01-29", 10, 0.068, "29.01.202610"),
  @("2026-01-29", 11, 0.382, "29.01.202611"),
  @("2026-01-29", 12, 0.576, "29.01.202612"),
  @("2026-01-29", 13, 0.878, "29.01.202613"),
  @("2026-01-29", 14, 0.872, "29.01.202614"),
  @("2026-01-29", 15, 0.603, "29.01.202615"),
  @("2026-01-29", 16, 0.391, "29.01.202616"),
  @("2026-01-29", 17, 0.227, "29.01.202617"),
  @("2026-01-29", 18, 0.058, "29.01.202618"),
  @("2026-01-29", 19, 0, "29.01.202619"),
  @("2026-01-29", 20, 0, "29.01.202620"),
  @("2026-01-29", 21, 0, "29.01.202621"),
  @("2026-01-29", 22, 0, "29.01.202622"),
  @("2026-01-29", 23, 0, "29.01.202623"),
  @("2026-01-29", 24, 0, "29.01.202624"),
  @("2026-01-30", 1, 0, "30.01.20261"),
  @("2026-01-30", 2, 0, "30.01.20262"),
  @("2026-01-30", 3, 0, "30.01.20263"),
  @("2026-01-30", 4, 0, "30.01.20264"),
  @("2026-01-30", 5, 0, "30.01.20265"),
  @("2026-01-30", 6, 0, "30.01.20266"),
  @("2026-01-30", 7, 0, "30.01.20267"),
  @("2026-01-30", 8, 0, "30.01.20268"),
  @("2026-01-30", 9, 0, "30.01.20269"),
  @("2026-01-30", 10, 0.048, "30.01.202610"),
  @("2026-01-30", 11, 0.261, "30.01.202611"),
  @("2026-01-30", 12, 0.498, "30.01.202612"),
  @("2026-01-30", 13, 0.885, "30.01.202613"),
  @("2026-01-30", 14, 0.972, "30.01.202614"),
  @("2026-01-30", 15, 0.906, "30.01.202615"),
  @("2026-01-30", 16, 0.528, "30.01.202616"),
  @("2026-01-30", 17, 0.302, "30.01.202617"),
  @("2026-01-30", 18, 0.097, "30.01.202618"),
  @("2026-01-30", 19, 0, "30.01.202619"),
  @("2026-01-30", 20, 0, "30.01.202620"),
  @("2026-01-30", 21, 0, "30.01.202621"),
  @("2026-01-30", 22, 0, "30.01.202622"),
  @("2026-01-30", 23, 0, "30.01.202623"),
  @("2026-01-30", 24, 0, "30.01.202624"),
  @("2026-01-31", 1, 0, "31.01.20261"),
  @("2026-01-31", 2, 0, "31.01.20262"),
  @("2026-01-31", 3, 0, "31.01.20263"),
  @("2026-01-31", 4, 0, "31.01.20264"),
  @("2026-01-31", 5, 0, "31.01.20265"),
  @("2026-01-31", 6, 0, "31.01.20266"),
  @("2026-01-31", 7, 0, "31.01.20267"),
  @("2026-01-31", 8, 0, "31.01.20268"),
  @("2026-01-31", 9, 0, "31.01.20269")
)

$nRows = $rows.Count
$nCols = 4
$data = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
  for ($j = 0; $j -lt $nCols; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}

$startRow = 2
$ws.Range("A$startRow").Resize($nRows, $nCols).Value = $data
